# Update "想去人数" (want-to-go headcount) figures in column F across the
# three sheets that carry this data: 展览, 演出, and 全部类型 (the latter
# mirrors the other sheets' rows, appended together).

$wb = $excel.ActiveWorkbook

# -- 展览 (sheet1) ----------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 303
$ws1.Range("F4").Value = 10288
$ws1.Range("F6").Value = 937
$ws1.Range("F7").Value = 1276
$ws1.Range("F8").Value = 6701
$ws1.Range("F9").Value = 15
$ws1.Range("F10").Value = 435
$ws1.Range("F11").Value = 197
$ws1.Range("F13").Value = 3163
$ws1.Range("F15").Value = 310
$ws1.Range("F16").Value = 640
$ws1.Range("F17").Value = 121
$ws1.Range("F18").Value = 689
$ws1.Range("F19").Value = 275
$ws1.Range("F20").Value = 57
$ws1.Range("F21").Value = 1612

# -- 演出 (sheet2) -----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28

# -- 全部类型 (sheet4) -------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 303
$ws4.Range("F4").Value = 10288
$ws4.Range("F6").Value = 937
$ws4.Range("F7").Value = 1276
$ws4.Range("F8").Value = 6703
$ws4.Range("F9").Value = 15
$ws4.Range("F10").Value = 435
$ws4.Range("F11").Value = 197
$ws4.Range("F13").Value = 3163
$ws4.Range("F15").Value = 310
$ws4.Range("F16").Value = 640
$ws4.Range("F17").Value = 121
$ws4.Range("F18").Value = 690
$ws4.Range("F19").Value = 275
$ws4.Range("F20").Value = 57
$ws4.Range("F21").Value = 1612
$ws4.Range("F22").Value = 28
